# update .bat and printing query dataframe
#
# The source dataframe (rows 2-9) was already duplicated once into rows
# 10-16. This edit appends one more duplicate of that same 7-row block
# (A10:AJ16) under it, as rows 17-23, growing the sheet's used range from
# A1:AJ16 to A1:AJ23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("A10:AJ16")
$destination = $ws.Range("A17:AJ23")

$source.Copy($destination)
